# Update the "dSF" (column F) values for the rows whose underlying data was
# re-pulled / recalculated. Column A holds a 0-based index that equals the
# worksheet row number minus 2 (data starts on row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row -> new value for column F (dSF)
$updates = @{
    2  = 0
    4  = 3
    7  = 6
    9  = -3
    10 = -4
    11 = -8
    12 = -3
    13 = 1
    14 = 3
    17 = 1
    20 = 2
    23 = 1
    24 = 1
    26 = -5
    29 = -2
    31 = 4
    35 = 4
    40 = 3
    45 = 1
    50 = -1
    51 = 3
    54 = 0
    56 = -10
    65 = 8
    68 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
